# Updated cryptos list on Sun Dec 10 09:37:07 UTC 2023 with GitHub Actions
# Refresh price / 1h-volume-change figures (and two pairs of rows whose
# ranking swapped places) on the "cryptos" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.659.68'
$ws.Range("E2").Value = '  -0.71%  '
$ws.Range("D3").Value = '2.337.88'
$ws.Range("E3").Value = '  -1.55%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '239.02'
$ws.Range("E5").Value = '  -1.44%  '
$ws.Range("D6").Value = '0.663'
$ws.Range("E6").Value = '  -4.02%  '
$ws.Range("D7").Value = '71.85'
$ws.Range("E7").Value = '  -6.26%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.586'
$ws.Range("E9").Value = '  -8.22%  '
$ws.Range("D10").Value = '0.0987'
$ws.Range("E10").Value = '  -4.39%  '
$ws.Range("D11").Value = '58.29'
$ws.Range("E11").Value = '  +1.28%  '
$ws.Range("D12").Value = '32.25'
$ws.Range("E12").Value = '  -2.75%  '
$ws.Range("E13").Value = '  -0.82%  '
$ws.Range("D14").Value = '7.10'
$ws.Range("E14").Value = '  -7.12%  '
$ws.Range("D15").Value = '2.687.84'
$ws.Range("E15").Value = '  -1.48%  '
$ws.Range("D16").Value = '16.13'
$ws.Range("E16").Value = '  -5.11%  '
$ws.Range("D17").Value = '0.894'
$ws.Range("E17").Value = '  -3.98%  '
$ws.Range("D18").Value = '2.341.30'
$ws.Range("E18").Value = '  -1.25%  '
$ws.Range("D19").Value = '43.646.35'
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("E20").Value = '  -2.65%  '
$ws.Range("D21").Value = '6.64'
$ws.Range("E21").Value = '  -1.20%  '
$ws.Range("D22").Value = '77.81'
$ws.Range("E22").Value = '  -0.36%  '
$ws.Range("D23").Value = '251.52'
$ws.Range("E23").Value = '  -2.82%  '

# Rows 24/25 swapped places: ImmutableX <-> Dai
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("B25").Value = 'ImmutableX'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D25").Value = '1.90'
$ws.Range("E25").Value = '  +7.50%  '

$ws.Range("E26").Value = '  +1.34%  '
$ws.Range("E27").Value = '  -2.28%  '
$ws.Range("D28").Value = '10.30'
$ws.Range("E28").Value = '  -8.56%  '
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("D30").Value = '175.50'
$ws.Range("E30").Value = '  -0.32%  '
$ws.Range("D31").Value = '22.11'
$ws.Range("E31").Value = '  -4.93%  '
$ws.Range("E32").Value = '  -2.70%  '
$ws.Range("D33").Value = '0.134'
$ws.Range("E33").Value = '  -0.62%  '
$ws.Range("D34").Value = '0.0731'
$ws.Range("E34").Value = '  -3.29%  '
$ws.Range("D35").Value = '5.04'
$ws.Range("E35").Value = '  -5.92%  '
$ws.Range("D36").Value = '5.33'
$ws.Range("E36").Value = '  -1.02%  '
$ws.Range("D37").Value = '3.72'
$ws.Range("E37").Value = '  -3.15%  '

# Rows 38/39 swapped places: LidoDAOToken <-> THORChain
$ws.Range("B38").Value = 'THORChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D38").Value = '6.35'
$ws.Range("E38").Value = '  -2.86%  '
$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").Value = '2.36'
$ws.Range("E39").Value = '  -4.42%  '

$ws.Range("D40").Value = '0.0269'
$ws.Range("E40").Value = '  -2.69%  '
$ws.Range("D41").Value = '5.24'
$ws.Range("E41").Value = '  +16.15%  '
$ws.Range("D42").Value = '64.64'
$ws.Range("E42").Value = '  +17.50%  '
$ws.Range("D43").Value = '9.20'
$ws.Range("E43").Value = '  +1.95%  '
$ws.Range("E44").Value = '  +5.37%  '
$ws.Range("D45").Value = '18.72'
$ws.Range("E45").Value = '  -2.08%  '
$ws.Range("D46").Value = '0.195'
$ws.Range("E46").Value = '  -4.93%  '
$ws.Range("E47").Value = '  +0.09%  '

# Rows 48/49 swapped places: NEARProtocol <-> TrustWalletToken
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").Value = '1.22'
$ws.Range("E48").Value = '  -4.12%  '
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '2.43'
$ws.Range("E49").Value = '  -4.05%  '

$ws.Range("D50").Value = '1.14'
$ws.Range("E50").Value = '  -5.78%  '

# Row 51: Aave replaced by HuobiToken
$ws.Range("B51").Value = 'HuobiToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D51").Value = '2.90'
$ws.Range("E51").Value = '  +3.54%  '
